$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (stored widths differ from COM ColumnWidth by ~0.8333 chars padding)
$ws.Columns.Item(2).ColumnWidth = 6.166666666666667
$ws.Columns.Item(7).ColumnWidth = 6.166666666666667
$ws.Columns.Item(13).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 6.166666666666667
$ws.Columns.Item(24).ColumnWidth = 6.166666666666667

# Update data rows 2-5 with new sensor readings
$ws.Cells.Item(2,1).Value = 45109.50694444445
$ws.Cells.Item(2,2).Value = 12.77
$ws.Cells.Item(2,3).Value = 8.425000000000001
$ws.Cells.Item(2,4).Value = 3.488
$ws.Cells.Item(2,5).Value = 27.776
$ws.Cells.Item(2,6).Value = 20.765
$ws.Cells.Item(2,7).Value = 9.82
$ws.Cells.Item(2,8).Value = 29.191
$ws.Cells.Item(2,9).Value = 15.707
$ws.Cells.Item(2,10).Value = 6.218
$ws.Cells.Item(2,11).Value = 9.112
$ws.Cells.Item(2,12).Value = 10.928
$ws.Cells.Item(2,13).Value = 11.76
$ws.Cells.Item(2,14).Value = 3.255
$ws.Cells.Item(2,15).Value = 10.151
$ws.Cells.Item(2,16).Value = 13.871
$ws.Cells.Item(2,17).Value = 9.122999999999999
$ws.Cells.Item(2,18).Value = 2.77
$ws.Cells.Item(2,19).Value = 1.592
$ws.Cells.Item(2,20).Value = 147.125
$ws.Cells.Item(2,21).Value = 28.067
$ws.Cells.Item(2,22).Value = 9.369999999999999
$ws.Cells.Item(2,23).Value = 18.004
$ws.Cells.Item(2,24).Value = 9.226000000000001
$ws.Cells.Item(2,25).Value = 2.728
$ws.Cells.Item(2,26).Value = 15.858
$ws.Cells.Item(2,27).Value = 8.276
$ws.Cells.Item(2,28).Value = 7.652
$ws.Cells.Item(2,29).Value = 9.069000000000001
$ws.Cells.Item(2,30).Value = 11.378
$ws.Cells.Item(2,31).Value = 2.869
$ws.Cells.Item(2,32).Value = 26.432
$ws.Cells.Item(2,33).Value = 4.957
$ws.Cells.Item(2,34).Value = 11.713

$ws.Cells.Item(3,1).Value = 45109.51388888889
$ws.Cells.Item(3,2).Value = 9.454000000000001
$ws.Cells.Item(3,3).Value = 6.662
$ws.Cells.Item(3,4).Value = 1.528
$ws.Cells.Item(3,5).Value = 20.843
$ws.Cells.Item(3,6).Value = 16.056
$ws.Cells.Item(3,7).Value = 7.317
$ws.Cells.Item(3,8).Value = 29.07
$ws.Cells.Item(3,9).Value = 11.635
$ws.Cells.Item(3,10).Value = 4.876
$ws.Cells.Item(3,11).Value = 6.977
$ws.Cells.Item(3,12).Value = 8.298
$ws.Cells.Item(3,13).Value = 8.936999999999999
$ws.Cells.Item(3,14).Value = 2.417
$ws.Cells.Item(3,15).Value = 7.519
$ws.Cells.Item(3,16).Value = 10.463
$ws.Cells.Item(3,17).Value = 6.744
$ws.Cells.Item(3,18).Value = 1.322
$ws.Cells.Item(3,19).Value = 0.8090000000000001
$ws.Cells.Item(3,20).Value = 107.128
$ws.Cells.Item(3,21).Value = 21.08
$ws.Cells.Item(3,22).Value = 6.941
$ws.Cells.Item(3,23).Value = 13.716
$ws.Cells.Item(3,24).Value = 7.098
$ws.Cells.Item(3,25).Value = 1.576
$ws.Cells.Item(3,26).Value = 14.478
$ws.Cells.Item(3,27).Value = 6.131
$ws.Cells.Item(3,28).Value = 5.642
$ws.Cells.Item(3,29).Value = 6.645
$ws.Cells.Item(3,30).Value = 8.638
$ws.Cells.Item(3,31).Value = 1.165
$ws.Cells.Item(3,32).Value = 26.904
$ws.Cells.Item(3,33).Value = 3.714
$ws.Cells.Item(3,34).Value = 8.678000000000001

$ws.Cells.Item(4,1).Value = 45109.52083333334
$ws.Cells.Item(4,2).Value = 6.601
$ws.Cells.Item(4,3).Value = 4.716
$ws.Cells.Item(4,4).Value = 0.976
$ws.Cells.Item(4,5).Value = 14.614
$ws.Cells.Item(4,6).Value = 11.27
$ws.Cells.Item(4,7).Value = 5.115
$ws.Cells.Item(4,8).Value = 21.867
$ws.Cells.Item(4,9).Value = 8.144
$ws.Cells.Item(4,10).Value = 3.44
$ws.Cells.Item(4,11).Value = 4.88
$ws.Cells.Item(4,12).Value = 5.834
$ws.Cells.Item(4,13).Value = 6.302
$ws.Cells.Item(4,14).Value = 1.693
$ws.Cells.Item(4,15).Value = 5.263
$ws.Cells.Item(4,16).Value = 7.34
$ws.Cells.Item(4,17).Value = 4.736
$ws.Cells.Item(4,18).Value = 0.88
$ws.Cells.Item(4,19).Value = 0.524
$ws.Cells.Item(4,20).Value = 72.792
$ws.Cells.Item(4,21).Value = 14.792
$ws.Cells.Item(4,22).Value = 4.858
$ws.Cells.Item(4,23).Value = 9.632
$ws.Cells.Item(4,24).Value = 4.997
$ws.Cells.Item(4,25).Value = 1.081
$ws.Cells.Item(4,26).Value = 10.532
$ws.Cells.Item(4,27).Value = 4.291
$ws.Cells.Item(4,28).Value = 3.955
$ws.Cells.Item(4,29).Value = 4.646
$ws.Cells.Item(4,30).Value = 6.082
$ws.Cells.Item(4,31).Value = 0.733
$ws.Cells.Item(4,32).Value = 20.176
$ws.Cells.Item(4,33).Value = 2.592
$ws.Cells.Item(4,34).Value = 6.075

$ws.Cells.Item(5,1).Value = 45109.52777777778
$ws.Cells.Item(5,2).Value = 11.43
$ws.Cells.Item(5,3).Value = 8.42
$ws.Cells.Item(5,4).Value = 0.93
$ws.Cells.Item(5,5).Value = 25.08
$ws.Cells.Item(5,6).Value = 20.11
$ws.Cells.Item(5,7).Value = 8.93
$ws.Cells.Item(5,8).Value = 33.63
$ws.Cells.Item(5,9).Value = 13.96
$ws.Cells.Item(5,10).Value = 6.1
$ws.Cells.Item(5,11).Value = 8.92
$ws.Cells.Item(5,12).Value = 10.05
$ws.Cells.Item(5,13).Value = 10.74
$ws.Cells.Item(5,14).Value = 2.9
$ws.Cells.Item(5,15).Value = 9.02
$ws.Cells.Item(5,16).Value = 12.74
$ws.Cells.Item(5,17).Value = 7.77
$ws.Cells.Item(5,18).Value = 0.72
$ws.Cells.Item(5,19).Value = 0.57
$ws.Cells.Item(5,20).Value = 130
$ws.Cells.Item(5,21).Value = 25.17
$ws.Cells.Item(5,22).Value = 8.33
$ws.Cells.Item(5,23).Value = 16.75
$ws.Cells.Item(5,24).Value = 8.779999999999999
$ws.Cells.Item(5,25).Value = 1.48
$ws.Cells.Item(5,26).Value = 16.51
$ws.Cells.Item(5,27).Value = 7.36
$ws.Cells.Item(5,28).Value = 6.6
$ws.Cells.Item(5,29).Value = 7.76
$ws.Cells.Item(5,30).Value = 10.53
$ws.Cells.Item(5,31).Value = 0.54
$ws.Cells.Item(5,32).Value = 30.53
$ws.Cells.Item(5,33).Value = 4.61
$ws.Cells.Item(5,34).Value = 10.41

# Remove row 6 (dataset now has 4 data rows instead of 5)
$ws.Rows.Item(6).Delete()